$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 53 values that changed
$ws.Range("F53").Value = 73393
$ws.Range("H53").Value = 44148
$ws.Range("T53").Value = 177316
$ws.Range("W53").Value = 34096
$ws.Range("Y53").Value = 18094
$ws.Range("AA53").Value = 63109
$ws.Range("AB53").Value = 240426

# Add new row 54 with data for 01-04-2021
# The label looks like a date, so a plain .Value assignment would get
# auto-converted to a date serial by Excel. Enter it as a text formula first,
# then paste-special as values only so it lands as a literal text string
# (shared string), matching the existing "Serie" column entries, without
# picking up a date number-format style along the way.
$ws.Range("A54").Formula = "=""01-04-2021"""
$ws.Range("A54").Copy()
$ws.Range("A54").PasteSpecial(-4163)
$excel.CutCopyMode = $false
$ws.Range("B54").Value = 31868
$ws.Range("C54").Value = 31300
$ws.Range("D54").Value = 568
$ws.Range("E54").Value = 397
$ws.Range("F54").Value = 70629
$ws.Range("G54").Value = 26783
$ws.Range("H54").Value = 43846
$ws.Range("I54").Value = 7454
$ws.Range("J54").Value = 1740
$ws.Range("K54").Value = 483
$ws.Range("L54").Value = 1249
$ws.Range("M54").Value = 3982
$ws.Range("N54").Value = 24049
$ws.Range("O54").Value = 3410
$ws.Range("P54").Value = 20639
$ws.Range("Q54").Value = 51620
$ws.Range("R54").Value = 22938
$ws.Range("S54").Value = 149
$ws.Range("T54").Value = 185769
$ws.Range("U54").Value = 8996
$ws.Range("V54").Value = 3045
$ws.Range("W54").Value = 34584
$ws.Range("X54").Value = 16131
$ws.Range("Y54").Value = 18453
$ws.Range("Z54").Value = 20120
$ws.Range("AA54").Value = 66745
$ws.Range("AB54").Value = 252513
